# Skill.xlsx update — rename the "normal attack" / "skill" row labels to their
# upper-cased forms, repoint the NextLevelID column (C) at the upper-cased
# "NORMALATTACK2" label for every data row, bump the AnimaState (G) values
# for the newly added skill rows, and refresh the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: row labels (ID), upper-cased ------------------------------
$ws.Range("A2").Value = "NORMALATTACK1"
$ws.Range("A3").Value = "NORMALATTACK2"
$ws.Range("A4").Value = "NORMALATTACK3"
$ws.Range("A5").Value = "NORMALTHUMP"
$ws.Range("A6").Value = "SKILL1"
$ws.Range("A7").Value = "SKILL2"
$ws.Range("A8").Value = "SKILL3"
$ws.Range("A9").Value = "SKILL4"

# --- Column C: NextLevelID, now uniformly "NORMALATTACK2" ----------------
# Clear the stray highlight fill/border that a few of these cells carried
# (rows 4-9) and make sure the column keeps its Text number format.
$ws.Range("C2:C9").ClearFormats()
$ws.Range("C2:C9").Value = "NORMALATTACK2"
$ws.Range("C2:C9").NumberFormat = "@"

# --- Column G: AnimaState bumped by 1 for the new skill rows -------------
$ws.Range("G6").Value = 101
$ws.Range("G7").Value = 102
$ws.Range("G8").Value = 103
$ws.Range("G9").Value = 104

# --- Selection: the user last worked on the NextLevelID column -----------
$null = $ws.Range("C2:C9").Select()
